$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("A18").Value = 20
$ws.Range("B18").Value = 63.56
$ws.Range("C18").Value = 12
$ws.Range("D18").Value = 300
$ws.Range("E18").Value = 8000
$ws.Range("F18").Value = "OPEXUM GCKŁNŻ AĄBĆDĘ FHIJLŃ ÓQRSŚT VWYZŹ|"
$ws.Range("G18").Value = -1403.88365943578
$ws.Range("H18").Value = "OPEXUM GCKŁNŻ AĄBĆDĘ FHIJLŃ ÓQRSŚT VWYZŹ|"
$ws.Range("I18").Value = -1403.8837

# Row 19
$ws.Range("A19").Value = 24
$ws.Range("B19").Value = 82.25
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 300
$ws.Range("E19").Value = 8000
$ws.Range("F19").Value = "HGZŁJU BSŹAT| ĄCĆDEĘ FIKLMN ŃOÓPQR ŚVWXYŻ"
$ws.Range("G19").Value = -1403.88365943578
$ws.Range("H19").Value = "HGZŁJU BSŹAT| ĄCĆDEĘ FIKLMN ŃOÓPQR ŚVWXYŻ"
$ws.Range("I19").Value = -1403.8837

# Row 20
$ws.Range("A20").Value = 28
$ws.Range("B20").Value = 94.83
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 300
$ws.Range("E20").Value = 8000
$ws.Range("F20").Value = "TDŁFIS NHQ|UK AĄBCĆE ĘGJLMŃ OÓPRŚV WXYZŹŻ"
$ws.Range("G20").Value = -1403.88365943578
$ws.Range("H20").Value = "TDŁFIS NHQ|UK AĄBCĆE ĘGJLMŃ OÓPRŚV WXYZŹŻ"
$ws.Range("I20").Value = -1403.8837

# Row 21
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = 28.69
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 300
$ws.Range("E21").Value = 2000
$ws.Range("F21").Value = "VĄŁWŹY LHMIST ABCĆDE ĘFGJKN ŃOÓPQR ŚUXZŻ|"
$ws.Range("G21").Value = -1403.88365943578
$ws.Range("H21").Value = "VĄŁWŹY LHMIST ABCĆDE ĘFGJKN ŃOÓPQR ŚUXZŻ|"
$ws.Range("I21").Value = -1403.8837

# Row 22
$ws.Range("A22").Value = 35
$ws.Range("B22").Value = 72.20999999999999
$ws.Range("C22").Value = 12
$ws.Range("D22").Value = 300
$ws.Range("E22").Value = 2000
$ws.Range("F22").Value = "ÓIGUŹC DWĆASK ĄBEĘFH JLŁMNŃ OPQRŚT VXYZŻ|"
$ws.Range("G22").Value = -1403.88365943578
$ws.Range("H22").Value = "ÓIGUŹC DWĆASK ĄBEĘFH JLŁMNŃ OPQRŚT VXYZŻ|"
$ws.Range("I22").Value = -1403.8837

# Row 23
$ws.Range("A23").Value = 68
$ws.Range("B23").Value = 150.34
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 300
$ws.Range("E23").Value = 2000
$ws.Range("F23").Value = "ŹĆÓZMI FNWĄX| ABCDEĘ GHJKLŁ ŃOPQRS ŚTUVYŻ"
$ws.Range("G23").Value = -1403.88365943578
$ws.Range("H23").Value = "ŹĆÓZMI FNWĄX| ABCDEĘ GHJKLŁ ŃOPQRS ŚTUVYŻ"
$ws.Range("I23").Value = -1403.8837

# Row 24
$ws.Range("A24").Value = 20
$ws.Range("B24").Value = 47.52
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 300
$ws.Range("E24").Value = 2000
$ws.Range("F24").Value = "AGQDTĘ ŻĆŃC|V ĄBEFHI JKLŁMN OÓPRSŚ UWXYZŹ"
$ws.Range("G24").Value = -1403.88365943578
$ws.Range("H24").Value = "AGQDTĘ ŻĆŃC|V ĄBEFHI JKLŁMN OÓPRSŚ UWXYZŹ"
$ws.Range("I24").Value = -1403.8837

# Row 25
$ws.Range("A25").Value = 60
$ws.Range("B25").Value = 130.68
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 300
$ws.Range("E25").Value = 2000
$ws.Range("F25").Value = "N|ŹKŚF JŃMĘCR AĄBĆDE GHILŁO ÓPQSTU VWXYZŻ"
$ws.Range("G25").Value = -1403.88365943578
$ws.Range("H25").Value = "N|ŹKŚF JŃMĘCR AĄBĆDE GHILŁO ÓPQSTU VWXYZŻ"
$ws.Range("I25").Value = -1403.8837

# Row 26
$ws.Range("A26").Value = 147
$ws.Range("B26").Value = 328.28
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 300
$ws.Range("E26").Value = 2000
$ws.Range("F26").Value = "ZŹTÓŻQ ĄŚWSOĘ ABCĆDE FGHIJK LŁMNŃP RUVXY|"
$ws.Range("G26").Value = -1403.88365943578
$ws.Range("H26").Value = "ZŹTÓŻQ ĄŚWSOĘ ABCĆDE FGHIJK LŁMNŃP RUVXY|"
$ws.Range("I26").Value = -1403.8837

# Row 27
$ws.Range("A27").Value = 36
$ws.Range("B27").Value = 74.14
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 300
$ws.Range("E27").Value = 2000
$ws.Range("F27").Value = "MBNOXS FLVŃGĘ AĄCĆDE HIJKŁÓ PQRŚTU WYZŹŻ|"
$ws.Range("G27").Value = -1403.88365943578
$ws.Range("H27").Value = "MBNOXS FLVŃGĘ AĄCĆDE HIJKŁÓ PQRŚTU WYZŹŻ|"
$ws.Range("I27").Value = -1403.8837

# Row 28
$ws.Range("A28").Value = 52
$ws.Range("B28").Value = 106.34
$ws.Range("C28").Value = 12
$ws.Range("D28").Value = 300
$ws.Range("E28").Value = 2000
$ws.Range("F28").Value = "ZŃÓTLA OFŚNVM ĄBCĆDE ĘGHIJK ŁPQRSU WXYŹŻ|"
$ws.Range("G28").Value = -1403.88365943578
$ws.Range("H28").Value = "ZŃÓTLA OFŚNVM ĄBCĆDE ĘGHIJK ŁPQRSU WXYŹŻ|"
$ws.Range("I28").Value = -1403.8837

# Row 29
$ws.Range("A29").Value = 33
$ws.Range("B29").Value = 70.59
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 300
$ws.Range("E29").Value = 2000
$ws.Range("F29").Value = "OŹĘŁGZ LMPHŚŻ AĄBCĆD EFIJKN ŃÓQRST UVWXY|"
$ws.Range("G29").Value = -1403.88365943578
$ws.Range("H29").Value = "OŹĘŁGZ LMPHŚŻ AĄBCĆD EFIJKN ŃÓQRST UVWXY|"
$ws.Range("I29").Value = -1403.8837

# Row 30
$ws.Range("A30").Value = 68
$ws.Range("B30").Value = 125.79
$ws.Range("C30").Value = 12
$ws.Range("D30").Value = 300
$ws.Range("E30").Value = 2000
$ws.Range("F30").Value = "AWCJQN KUPYIV ĄBĆDEĘ FGHLŁM ŃOÓRSŚ TXZŹŻ|"
$ws.Range("G30").Value = -1403.88365943578
$ws.Range("H30").Value = "AWCJQN KUPYIV ĄBĆDEĘ FGHLŁM ŃOÓRSŚ TXZŹŻ|"
$ws.Range("I30").Value = -1403.8837

# Row 31
$ws.Range("A31").Value = 274
$ws.Range("B31").Value = "Attempt failed!"

# Row 32
$ws.Range("A32").Value = 274
$ws.Range("B32").Value = 643.91
$ws.Range("C32").Value = 12
$ws.Range("D32").Value = 300
$ws.Range("E32").Value = 2000
$ws.Range("F32").Value = "AZDŃEQ YLKMPN ĄBCĆĘF GHIJŁO ÓRSŚTU VWXŹŻ|"
$ws.Range("G32").Value = -1403.88365943578
$ws.Range("H32").Value = "VBJFTŻ UÓRWZS Ź|OMGQ NĄKŚLŁ CAHXDE ŃĘPĆIY"
$ws.Range("I32").Value = -2001.6128

# Row 33
$ws.Range("A33").Value = 20
$ws.Range("B33").Value = 46.08
$ws.Range("C33").Value = 12
$ws.Range("D33").Value = 300
$ws.Range("E33").Value = 2000
$ws.Range("F33").Value = "ŚJŁĆDŹ ĘŻHVFQ AĄBCEG IKLMNŃ OÓPRST UWXYZ|"
$ws.Range("G33").Value = -1403.88365943578
$ws.Range("H33").Value = "ŚJŁĆDV ĘŻHŹFQ AĄBCEG IKLMNŃ OÓPRST UWXYZ|"
$ws.Range("I33").Value = -1437.7579

# Row 34
$ws.Range("A34").Value = 37
$ws.Range("B34").Value = 80.31
$ws.Range("C34").Value = 12
$ws.Range("D34").Value = 300
$ws.Range("E34").Value = 2000
$ws.Range("F34").Value = "QONLĆH JCÓSPV AĄBDEĘ FGIKŁM ŃRŚTUW XYZŹŻ|"
$ws.Range("G34").Value = -1403.88365943578
$ws.Range("H34").Value = "FGIKŁM ŃRŚTUW XYZŹŻ| QONLĆH JCÓSPV AĄBDEĘ"
$ws.Range("I34").Value = -1403.8837

# Row 35
$ws.Range("A35").Value = 255
$ws.Range("B35").Value = "Attempt failed!"

# Row 36
$ws.Range("A36").Value = 255
$ws.Range("B36").Value = 629.58
$ws.Range("C36").Value = 12
$ws.Range("D36").Value = 300
$ws.Range("E36").Value = 2000
$ws.Range("F36").Value = "DKORCŁ I|ĘLHA ĄBĆEFG JMNŃÓP QSŚTUV WXYZŹŻ"
$ws.Range("G36").Value = -1403.88365943578
$ws.Range("H36").Value = "ZLTROŃ JCWPŹÓ YGXKŻF QĆMŚBE ADĘNĄH |IVSŁU"
$ws.Range("I36").Value = -1976.7157

# Row 37
$ws.Range("A37").Value = 27
$ws.Range("B37").Value = 58.3
$ws.Range("C37").Value = 12
$ws.Range("D37").Value = 300
$ws.Range("E37").Value = 2000
$ws.Range("F37").Value = "WŚOS|Ę DKVÓCĆ AĄBEFG HIJLŁM NŃPQRT UXYZŹŻ"
$ws.Range("G37").Value = -1403.88365943578
$ws.Range("H37").Value = "WŚOS|Ę DKVÓCĆ AĄBEFG HIJLŁM NŃPQRT UXYZŹŻ"
$ws.Range("I37").Value = -1403.8837
